# Update the two-digit multiplication equations throughout the document.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "40×56="; new = "59×15="},
    @{old = "55×14="; new = "38×55="},
    @{old = "46×26="; new = "35×20="},
    @{old = "78×42="; new = "57×81="},
    @{old = "11×45="; new = "58×99="},
    @{old = "39×91="; new = "19×23="},
    @{old = "78×15="; new = "50×80="},
    @{old = "36×72="; new = "58×37="},
    @{old = "45×38="; new = "18×34="},
    @{old = "46×86="; new = "31×52="},
    @{old = "14×70="; new = "97×62="},
    @{old = "17×70="; new = "64×51="},
    @{old = "22×84="; new = "97×92="},
    @{old = "46×93="; new = "63×47="},
    @{old = "83×59="; new = "94×49="},
    @{old = "56×13="; new = "94×67="},
    @{old = "65×18="; new = "96×12="},
    @{old = "91×31="; new = "73×76="},
    @{old = "21×87="; new = "43×80="},
    @{old = "48×30="; new = "49×16="},
    @{old = "90×13="; new = "93×15="},
    @{old = "12×40="; new = "58×94="},
    @{old = "90×97="; new = "94×60="},
    @{old = "73×31="; new = "33×57="},
    @{old = "25×87="; new = "90×14="}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
